$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are stored as plain text in the source sheet (some rows use
# "." as a thousands separator rather than a decimal point). Force each Price cell that
# is being updated to Text format first so Excel does not silently convert a value such
# as "1.002" into the number 1.002, which would lose the original text formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.189.61'
$ws.Range("E2").Value = '  -0.77%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.676.54'
$ws.Range("E3").Value = '  -0.98%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.55'
$ws.Range("E5").Value = '  -3.42%  '

$ws.Range("E6").Value = '  -3.93%  '

$ws.Range("E7").Value = '  -0.81%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2662'
$ws.Range("E8").Value = '  -3.26%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06313'
$ws.Range("E9").Value = '  -2.21%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.38'
$ws.Range("E10").Value = '  -2.86%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07535'
$ws.Range("E11").Value = '  -1.93%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.680.31'
$ws.Range("E12").Value = '  -0.67%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.458'
$ws.Range("E13").Value = '  -1.85%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5655'
$ws.Range("E14").Value = '  -3.16%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008052'
$ws.Range("E15").Value = '  -3.66%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.57'
$ws.Range("E16").Value = '  +1.68%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.256.30'
$ws.Range("E17").Value = '  -0.66%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.002'
$ws.Range("E18").Value = '  -0.82%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.825'
$ws.Range("E19").Value = '  -2.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '188.54'
$ws.Range("E20").Value = '  -1.72%  '

$ws.Range("E21").Value = '  -5.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.208'
$ws.Range("E22").Value = '  -0.71%  '

$ws.Range("E23").Value = '  -0.83%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '148.48'
$ws.Range("E24").Value = '  -0.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1255'
$ws.Range("E25").Value = '  -5.22%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.599'
$ws.Range("E26").Value = '  -4.03%  '

$ws.Range("E27").Value = '  +1.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06238'
$ws.Range("E28").Value = '  -0.85%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.358'
$ws.Range("E29").Value = '  -1.89%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.282'
$ws.Range("E30").Value = '  -3.72%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.497'
$ws.Range("E31").Value = '  -3.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.448'
$ws.Range("E32").Value = '  -4.39%  '

$ws.Range("E33").Value = '  -3.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.003'
$ws.Range("E34").Value = '  -3.87%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6069'
$ws.Range("E35").Value = '  -1.37%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.403'
$ws.Range("E36").Value = '  -0.49%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.723'
$ws.Range("E37").Value = '  +0.38%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.123'
$ws.Range("E38").Value = '  -1.10%  '

$ws.Range("E39").Value = '  -1.56%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.081.64'
$ws.Range("E40").Value = '  -3.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8685'
$ws.Range("E41").Value = '  -1.53%  '

$ws.Range("E42").Value = '  -1.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.15'
$ws.Range("E43").Value = '  -1.38%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.826.49'
$ws.Range("E44").Value = '  -0.97%  '

$ws.Range("E45").Value = '  +1.34%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.27'
$ws.Range("E46").Value = '  -2.22%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  -0.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.045'
$ws.Range("E48").Value = '  -1.84%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05240'
$ws.Range("E49").Value = '  -0.63%  '

$ws.Range("E50").Value = '  -1.19%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.990'
$ws.Range("E51").Value = '  -1.90%  '

